$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns per latest crypto data pull ---
$ws.Range("D2").Value = "63.243.75"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.567.43"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.17"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.05"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("E9").Value = "  +3.54%  "
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.93"
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("D14").Value = "3.028.67"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "63.140.01"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("E16").Value = "  +3.78%  "
$ws.Range("D17").Value = "2.570.26"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "341.49"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("E20").Value = "  +2.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.88"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.23"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("E24").Value = "  +4.30%  "
$ws.Range("D25").Value = "2.687.82"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.27"
$ws.Range("E27").Value = "  +15.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.58"
$ws.Range("E28").Value = "  +1.94%  "

# Rows 29 and 30 swapped position in the upstream ranking (Binance-PegBSC-USD <-> SuiNetwork)
$ws.Range("B29").Value = "SuiNetwork"
$ws.Range("C29").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.49"
$ws.Range("E29").Value = "  +1.22%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("E31").Value = "  +5.46%  "
$ws.Range("D32").Value = "0.0₃0831"
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "177.66"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "442.62"
$ws.Range("E34").Value = "  +6.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.409"
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.39"
$ws.Range("E37").Value = "  +2.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.53"
$ws.Range("E38").Value = "  +3.13%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "152.39"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.86"
$ws.Range("E43").Value = "  +3.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.54"
$ws.Range("E44").Value = "  +4.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0555"
$ws.Range("E45").Value = "  +6.26%  "
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.54"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("E51").Value = "  -0.31%  "
